$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.714.40'
$ws.Range('E2').Value = '  -2.46%  '
$ws.Range('D3').Value = '3.282.21'
$ws.Range('E3').Value = '  -0.98%  '
$ws.Range('E4').Value = '  +0.04%  '
$style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '573.05'
$ws.Range('D5').Style = $style
$ws.Range('E5').Value = '  -0.90%  '
$style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '177.78'
$ws.Range('D6').Style = $style
$ws.Range('E6').Value = '  -4.50%  '
$style = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.629'
$ws.Range('D7').Style = $style
$ws.Range('E7').Value = '  +4.32%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -2.91%  '
$ws.Range('E10').Value = '  +0.66%  '
$style = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.400'
$ws.Range('D11').Style = $style
$ws.Range('E11').Value = '  -2.81%  '
$ws.Range('D12').Value = '3.856.12'
$ws.Range('E12').Value = '  -0.89%  '
$ws.Range('E13').Value = '  -3.57%  '
$style = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.57'
$ws.Range('D14').Style = $style
$ws.Range('E14').Value = '  -3.41%  '
$ws.Range('D15').Value = '65.892.00'
$ws.Range('E15').Value = '  -2.54%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.301.02'
$ws.Range('E16').Value = '  -0.31%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$style = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000163'
$ws.Range('D17').Style = $style
$ws.Range('E17').Value = '  -2.63%  '
$style = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '436.28'
$ws.Range('D18').Style = $style
$ws.Range('E18').Value = '  -1.66%  '
$style = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.57'
$ws.Range('D19').Style = $style
$ws.Range('E19').Value = '  -2.71%  '
$style = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.20'
$ws.Range('D20').Style = $style
$ws.Range('E20').Value = '  -2.72%  '
$style = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.40'
$ws.Range('D21').Style = $style
$ws.Range('E21').Value = '  -4.88%  '
$style = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.41'
$ws.Range('D22').Style = $style
$ws.Range('E22').Value = '  -2.09%  '
$style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('D23').Style = $style
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').Value = '3.434.14'
$ws.Range('E24').Value = '  -0.66%  '
$style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.509'
$ws.Range('D25').Style = $style
$ws.Range('E25').Value = '  -0.98%  '
$ws.Range('E26').Value = '  -5.19%  '
$ws.Range('E27').Value = '  +2.96%  '
$style = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.89'
$ws.Range('D28').Style = $style
$ws.Range('E28').Value = '  -2.03%  '
$style = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = $style
$ws.Range('E29').Value = '  -0.01%  '
$style = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.93'
$ws.Range('D30').Style = $style
$ws.Range('E30').Value = '  -2.32%  '
$style = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.32'
$ws.Range('D31').Style = $style
$ws.Range('E31').Value = '  -2.70%  '
$ws.Range('E32').Value = '  +0.10%  '
$style = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.14'
$ws.Range('D33').Style = $style
$ws.Range('E33').Value = '  -3.92%  '
$style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.62'
$ws.Range('D34').Style = $style
$ws.Range('E34').Value = '  -2.99%  '
$ws.Range('E35').Value = '  -4.92%  '
$style = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '158.73'
$ws.Range('D36').Style = $style
$ws.Range('E36').Value = '  -2.53%  '
$ws.Range('E37').Value = '  -5.01%  '
$style = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '26.73'
$ws.Range('D38').Style = $style
$ws.Range('E38').Value = '  -2.04%  '
$style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.78'
$ws.Range('D39').Style = $style
$ws.Range('E39').Value = '  -4.01%  '
$ws.Range('D40').Value = '2.774.59'
$ws.Range('E40').Value = '  -0.29%  '
$style = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.780'
$ws.Range('D41').Style = $style
$ws.Range('E41').Value = '  -1.57%  '
$style = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.32'
$ws.Range('D42').Style = $style
$ws.Range('E42').Value = '  -3.74%  '
$style = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.32'
$ws.Range('D43').Style = $style
$ws.Range('E43').Value = '  +0.45%  '
$ws.Range('E44').Value = '  -3.82%  '
$style = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0656'
$ws.Range('D45').Style = $style
$ws.Range('E45').Value = '  -2.63%  '
$ws.Range('E46').Value = '  -4.88%  '
$style = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '320.64'
$ws.Range('D47').Style = $style
$ws.Range('E47').Value = '  -2.43%  '
$style = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.42'
$ws.Range('D48').Style = $style
$ws.Range('E48').Value = '  -6.00%  '
$ws.Range('E49').Value = '  -1.90%  '
$ws.Range('E50').Value = '  +2.25%  '
$ws.Range('E51').Value = '  +0.08%  '
